$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.089.87"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.267.78"
$ws.Range("E3").Value = "  -1.50%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.25"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.97"
$ws.Range("E6").Value = "  -4.45%  "

# Row 7
$ws.Range("E7").Value = "  -2.64%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -2.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.96"
$ws.Range("E10").Value = "  -4.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -0.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.39"
$ws.Range("E12").Value = "  -6.72%  "

# Row 13
$ws.Range("E13").Value = "  +0.93%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.78"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.621.19"
$ws.Range("E16").Value = "  -1.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.266.38"
$ws.Range("E17").Value = "  -0.66%  "

# Row 18
$ws.Range("E18").Value = "  -2.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.090.69"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  +1.78%  "

# Row 21
$ws.Range("E21").Value = "  -1.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.97"
$ws.Range("E22").Value = "  -1.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.16"
$ws.Range("E23").Value = "  -2.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.82"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("E25").Value = "  -0.79%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -2.70%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.71"
$ws.Range("E28").Value = "  -5.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  -2.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.59"
$ws.Range("E30").Value = "  +2.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.66"
$ws.Range("E31").Value = "  -3.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.14"
$ws.Range("E32").Value = "  -0.32%  "

# Row 33
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("E34").Value = "  +6.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.88"
$ws.Range("E35").Value = "  -2.60%  "

# Row 36
$ws.Range("E36").Value = "  -1.91%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.76"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0686"
$ws.Range("E38").Value = "  -2.56%  "

# Row 39
$ws.Range("E39").Value = "  -2.89%  "

# Row 40
$ws.Range("E40").Value = "  -1.58%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("E41").Value = "  -2.50%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.72"
$ws.Range("E42").Value = "  -4.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -7.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.956.70"
$ws.Range("E44").Value = "  -0.50%  "

# Row 45
$ws.Range("E45").Value = "  -1.42%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.55"
$ws.Range("E46").Value = "  -4.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.56"
$ws.Range("E47").Value = "  -6.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -4.23%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.493.46"
$ws.Range("E49").Value = "  -1.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.91"
$ws.Range("E50").Value = "  -6.34%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.73"
$ws.Range("E51").Value = "  -3.95%  "
